$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, shifting existing rows 100:216 down to 101:217
$ws.Rows("100:100").Insert()

# Populate the newly inserted row 100 with the new price report record
$ws.Range("A100").Value = 7
$ws.Range("B100").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C100").Value = "Ñuble"
$ws.Range("D100").Value = 44763
$ws.Range("E100").Value = 16
$ws.Range("F100").Value = 100112017
$ws.Range("G100").Value = "Apio"
$ws.Range("H100").Value = "Americana (o)"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 120
$ws.Range("K100").Value = 9000
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = 9500
$ws.Range("N100").Value = "`$/docena de matas"
$ws.Range("O100").Value = "Provincia del Elquí"
$ws.Range("P100").Value = 1583
$ws.Range("Q100").Value = 6
$ws.Range("R100").Value = "Hortaliza"
